# Auto-generated: apply value updates to match target diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 125711.375
$ws.Range("I33").Value = 143527.28
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 143527.28
$ws.Range("L33").Value = 1000
$ws.Range("M33").Value = -143298.28
$ws.Range("N33").Value = -1458

$ws.Range("H98").Value = 1266.6666
$ws.Range("I98").Value = 1295.6522
$ws.Range("J98").Value = 600
$ws.Range("K98").Value = 1295.6522
$ws.Range("L98").Value = 600
$ws.Range("M98").Value = 202.3478
$ws.Range("N98").Value = -3596

$ws.Range("H106").Value = 2418.3635
$ws.Range("I106").Value = 2432.889
$ws.Range("K106").Value = 2432.889
$ws.Range("M106").Value = -1801.889

$ws.Range("H122").Value = 1266.6666
$ws.Range("I122").Value = 1295.6522
$ws.Range("J122").Value = 600
$ws.Range("K122").Value = 3886.9566
$ws.Range("L122").Value = 1800
$ws.Range("M122").Value = -1436.9566
$ws.Range("N122").Value = -6700

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1590579.9
$ws.Range("I32").Value = 4504.1016
$ws.Range("J32").Value = 12534503
$ws.Range("K32").Value = 4504.1016
$ws.Range("L32").Value = 12534503
$ws.Range("M32").Value = -4217.1016
$ws.Range("N32").Value = -12535077

$ws.Range("H70").Value = 33499.5
$ws.Range("J70").Value = 33499.5
$ws.Range("L70").Value = 33499.5
$ws.Range("N70").Value = -34039.5

$ws.Range("H73").Value = 33499.5
$ws.Range("J73").Value = 33499.5
$ws.Range("L73").Value = 33499.5
$ws.Range("N73").Value = -35371.5

$ws.Range("H97").Value = 506.10715
$ws.Range("I97").Value = 498.95
$ws.Range("J97").Value = 524
$ws.Range("K97").Value = 498.95
$ws.Range("L97").Value = 524
$ws.Range("M97").Value = -2.949999999999989
$ws.Range("N97").Value = -1516

$ws.Range("H122").Value = 1276.591
$ws.Range("I122").Value = 1187.7333
$ws.Range("J122").Value = 1467
$ws.Range("K122").Value = 3563.199900000001
$ws.Range("L122").Value = 4401
$ws.Range("M122").Value = -1113.199900000001
$ws.Range("N122").Value = -9301

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 109702.9
$ws.Range("I134").Value = 5284.357
$ws.Range("J134").Value = 402074.8
$ws.Range("K134").Value = 15853.071
$ws.Range("L134").Value = 1206224.4
$ws.Range("M134").Value = -13318.071
$ws.Range("N134").Value = -1211294.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1114390.2
$ws.Range("I99").Value = 2502753
$ws.Range("K99").Value = 2502753
$ws.Range("M99").Value = -2501255

$ws.Range("H126").Value = 1114390.2
$ws.Range("I126").Value = 2502753
$ws.Range("K126").Value = 7508259
$ws.Range("M126").Value = -7505789

$ws.Range("H132").Value = 3608.7144
$ws.Range("I132").Value = 2673.4
$ws.Range("J132").Value = 4128.3335
$ws.Range("K132").Value = 8020.200000000001
$ws.Range("L132").Value = 12385.0005
$ws.Range("M132").Value = -5490.200000000001
$ws.Range("N132").Value = -17445.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 8783.333
$ws.Range("I44").Value = 300
$ws.Range("J44").Value = 10480
$ws.Range("K44").Value = 900
$ws.Range("L44").Value = 31440
$ws.Range("M44").Value = -502
$ws.Range("N44").Value = -32236

$ws.Range("H92").Value = 183
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 183
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 549
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -3045

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4625.375
$ws.Range("I70").Value = 4231.231
$ws.Range("J70").Value = 6333.3335
$ws.Range("K70").Value = 4231.231
$ws.Range("L70").Value = 6333.3335
$ws.Range("M70").Value = -3961.231
$ws.Range("N70").Value = -6873.3335

$ws.Range("H73").Value = 4625.375
$ws.Range("I73").Value = 4231.231
$ws.Range("J73").Value = 6333.3335
$ws.Range("K73").Value = 4231.231
$ws.Range("L73").Value = 6333.3335
$ws.Range("M73").Value = -3295.231
$ws.Range("N73").Value = -8205.3335

$ws.Range("H102").Value = 3495.8572
$ws.Range("J102").Value = 3599.6667
$ws.Range("L102").Value = 3599.6667
$ws.Range("N102").Value = -6843.6667

$ws.Range("H122").Value = 2319.9
$ws.Range("I122").Value = 1741.2941
$ws.Range("K122").Value = 5223.8823
$ws.Range("M122").Value = -2773.8823

$ws.Range("H126").Value = 7000
$ws.Range("I126").Value = 7000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 21000
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -18530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2446.8
$ws.Range("I7").Value = 2555.7778
$ws.Range("K7").Value = 2555.7778
$ws.Range("M7").Value = -2443.7778

$ws.Range("H40").Value = 3089.1
$ws.Range("I40").Value = 3284.4285
$ws.Range("K40").Value = 3284.4285
$ws.Range("M40").Value = -3148.4285

$ws.Range("H68").Value = 3912.75
$ws.Range("I68").Value = 3950.3333
$ws.Range("J68").Value = 3800
$ws.Range("K68").Value = 3950.3333
$ws.Range("L68").Value = 3800
$ws.Range("M68").Value = -3201.3333
$ws.Range("N68").Value = -5298

$ws.Range("H71").Value = 3912.75
$ws.Range("I71").Value = 3950.3333
$ws.Range("J71").Value = 3800
$ws.Range("K71").Value = 19751.6665
$ws.Range("L71").Value = 19000
$ws.Range("M71").Value = -16007.6665
$ws.Range("N71").Value = -26488

$ws.Range("H126").Value = 2446.8
$ws.Range("I126").Value = 2555.7778
$ws.Range("K126").Value = 7667.3334
$ws.Range("M126").Value = -5197.3334

$ws.Range("H132").Value = 3959.8667
$ws.Range("I132").Value = 3550
$ws.Range("J132").Value = 4428.2856
$ws.Range("K132").Value = 10650
$ws.Range("L132").Value = 13284.8568
$ws.Range("M132").Value = -8120
$ws.Range("N132").Value = -18344.8568

$ws.Range("H136").Value = 3528.842
$ws.Range("I136").Value = 1519.1428
$ws.Range("K136").Value = 4557.428400000001
$ws.Range("M136").Value = -2007.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1042.375
$ws.Range("I122").Value = 1075.8889
$ws.Range("J122").Value = 999.2857
$ws.Range("K122").Value = 3227.6667
$ws.Range("L122").Value = 2997.8571
$ws.Range("M122").Value = -777.6666999999998
$ws.Range("N122").Value = -7897.8571

$ws.Range("H132").Value = 3445.7273
$ws.Range("I132").Value = 3500
$ws.Range("J132").Value = 3400.5
$ws.Range("K132").Value = 10500
$ws.Range("L132").Value = 10201.5
$ws.Range("M132").Value = -7970
$ws.Range("N132").Value = -15261.5

$ws.Range("H136").Value = 1723.9286
$ws.Range("I136").Value = 1658.1482
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 4974.444600000001
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -2424.444600000001
$ws.Range("N136").Value = -15600

